$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 98.0026933829441
$ws.Range("C7").Value = 98.02712715519576
$ws.Range("D7").Value = 98.04793380265713
$ws.Range("E7").Value = 98.01521732164599

$ws.Range("B8").Value = 97.55810112812509
$ws.Range("C8").Value = 97.51237758389458
$ws.Range("D8").Value = 97.60548937726499
$ws.Range("E8").Value = 97.50472687945125

$ws.Range("B9").Value = 96.20948309419995
$ws.Range("C9").Value = 96.248303355997
$ws.Range("D9").Value = 96.20470324321657
$ws.Range("E9").Value = 96.21928597786848
